# Coursework2-part1.docx -- final proof-reading pass.
#
# The author re-opened the document, let Word's proofing pass run over it,
# and saved again ("Thank you, I am happy with this :)"). The wording of
# the document is unchanged; the only externally-visible effects are a
# couple of stale layout artifacts that get cleared out on resave:
#   - the floating "_GoBack" bookmark (left over from the last edit
#     session) disappears once that spot is touched again
#   - three stale <w:lastRenderedPageBreak/> markers (before "Plagiarism",
#     before "Value iteration computes ...", and before "Questions 5 & 6
#     ...") are dropped because the layout is recomputed
#
# Re-create both effects with no-op Find/Replace passes over the exact
# spots involved -- this forces Word to rewrite those runs (clearing the
# bookmark / stale page-break marker) while leaving the visible text
# identical.

$d = $word.ActiveDocument

# 1) Touch the text that straddles the old "_GoBack" bookmark so the
#    paragraph gets rewritten and the bookmark disappears.
$d.Content.Find.Execute(
    "comprised of the files below", $true, $false, $false, $false, $false,
    $true, 1, $false, "comprised of the files below", 1) | Out-Null

# 2) Re-touch the three headings/paragraphs that carried a stale
#    lastRenderedPageBreak marker so it gets cleared on re-layout.
$d.Content.Find.Execute(
    "Plagiarism", $true, $false, $false, $false, $false,
    $true, 1, $false, "Plagiarism", 1) | Out-Null

$d.Content.Find.Execute(
    "Value iteration computes k-step estimates of the optimal values, V",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Value iteration computes k-step estimates of the optimal values, V", 1) | Out-Null

$d.Content.Find.Execute(
    "Questions 5 & 6 are on Reinforcement Learning:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Questions 5 & 6 are on Reinforcement Learning:", 1) | Out-Null
